# Apply the Universalis/market-data refresh values captured in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 62
$ws.Range("H62").Value = 2476.923
$ws.Range("I62").Value = 2569.5652
$ws.Range("J62").Value = 1766.6666
$ws.Range("K62").Value = 2569.5652
$ws.Range("L62").Value = 1766.6666
$ws.Range("M62").Value = -1945.5652
$ws.Range("N62").Value = -3014.6666
# row 65
$ws.Range("H65").Value = 2476.923
$ws.Range("I65").Value = 2569.5652
$ws.Range("J65").Value = 1766.6666
$ws.Range("K65").Value = 12847.826
$ws.Range("L65").Value = 8833.333000000001
$ws.Range("M65").Value = -9727.826000000001
$ws.Range("N65").Value = -15073.333
# row 86
$ws.Range("H86").Value = 9581.385
$ws.Range("I86").Value = 3850.4443
$ws.Range("J86").Value = 22476
$ws.Range("K86").Value = 3850.4443
$ws.Range("L86").Value = 22476
$ws.Range("M86").Value = -2727.4443
$ws.Range("N86").Value = -24722
# row 89
$ws.Range("H89").Value = 9581.385
$ws.Range("I89").Value = 3850.4443
$ws.Range("J89").Value = 22476
$ws.Range("K89").Value = 19252.2215
$ws.Range("L89").Value = 112380
$ws.Range("M89").Value = -13636.2215
$ws.Range("N89").Value = -123612
# row 113
$ws.Range("H113").Value = 2849.8333
$ws.Range("I113").Value = 2866.3333
$ws.Range("J113").Value = 2833.3333
$ws.Range("K113").Value = 2866.3333
$ws.Range("L113").Value = 2833.3333
$ws.Range("M113").Value = 387.6667000000002
$ws.Range("N113").Value = -9341.3333
# row 129
$ws.Range("H129").Value = 132005.98
$ws.Range("I129").Value = 751296.5
$ws.Range("J129").Value = 1629.0264
$ws.Range("K129").Value = 2253889.5
$ws.Range("L129").Value = 4887.0792
$ws.Range("M129").Value = -2248889.5
$ws.Range("N129").Value = -14887.0792

$ws = $wb.Worksheets.Item("ARM")
# row 3
$ws.Range("H3").Value = 305
$ws.Range("I3").Value = 305
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 305
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -190
$ws.Range("N3").ClearContents()
# row 61
$ws.Range("H61").Value = 2758.6758
$ws.Range("I61").Value = 2595.48
$ws.Range("J61").Value = 3098.6667
$ws.Range("K61").Value = 2595.48
$ws.Range("L61").Value = 3098.6667
$ws.Range("M61").Value = -2383.48
$ws.Range("N61").Value = -3522.6667
# row 63
$ws.Range("H63").Value = 2865.4736
$ws.Range("I63").Value = 2442.9333
$ws.Range("J63").Value = 4450
$ws.Range("K63").Value = 2442.9333
$ws.Range("L63").Value = 4450
$ws.Range("M63").Value = -1756.9333
$ws.Range("N63").Value = -5822
# row 66
$ws.Range("H66").Value = 2865.4736
$ws.Range("I66").Value = 2442.9333
$ws.Range("J66").Value = 4450
$ws.Range("K66").Value = 12214.6665
$ws.Range("L66").Value = 22250
$ws.Range("M66").Value = -8782.666500000001
$ws.Range("N66").Value = -29114
# row 136
$ws.Range("H136").Value = 2758.6758
$ws.Range("I136").Value = 2595.48
$ws.Range("J136").Value = 3098.6667
$ws.Range("K136").Value = 7786.440000000001
$ws.Range("L136").Value = 9296.000100000001
$ws.Range("M136").Value = -5236.440000000001
$ws.Range("N136").Value = -14396.0001

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 4968.183
$ws.Range("I31").Value = 2140.3125
$ws.Range("J31").Value = 5996.5
$ws.Range("K31").Value = 2140.3125
$ws.Range("L31").Value = 5996.5
$ws.Range("M31").Value = -1845.3125
$ws.Range("N31").Value = -6586.5
# row 34
$ws.Range("H34").Value = 4968.183
$ws.Range("I34").Value = 2140.3125
$ws.Range("J34").Value = 5996.5
$ws.Range("K34").Value = 2140.3125
$ws.Range("L34").Value = 5996.5
$ws.Range("M34").Value = -1938.3125
$ws.Range("N34").Value = -6400.5
# row 58
$ws.Range("H58").Value = 2018.1897
$ws.Range("I58").Value = 1771.5106
$ws.Range("J58").Value = 3072.182
$ws.Range("K58").Value = 1771.5106
$ws.Range("L58").Value = 3072.182
$ws.Range("M58").Value = -1568.5106
$ws.Range("N58").Value = -3478.182
# row 99
$ws.Range("H99").Value = 2233.9048
$ws.Range("I99").Value = 2273.1428
$ws.Range("J99").Value = 2214.2856
$ws.Range("K99").Value = 2273.1428
$ws.Range("L99").Value = 2214.2856
$ws.Range("M99").Value = -775.1428000000001
$ws.Range("N99").Value = -5210.2856
# row 107
$ws.Range("H107").Value = 669.61536
$ws.Range("I107").Value = 587
$ws.Range("J107").Value = 1016.6
$ws.Range("K107").Value = 587
$ws.Range("L107").Value = 1016.6
$ws.Range("M107").Value = 1333
$ws.Range("N107").Value = -4856.6
# row 126
$ws.Range("H126").Value = 2233.9048
$ws.Range("I126").Value = 2273.1428
$ws.Range("J126").Value = 2214.2856
$ws.Range("K126").Value = 6819.428400000001
$ws.Range("L126").Value = 6642.8568
$ws.Range("M126").Value = -4349.428400000001
$ws.Range("N126").Value = -11582.8568
# row 136
$ws.Range("H136").Value = 2018.1897
$ws.Range("I136").Value = 1771.5106
$ws.Range("J136").Value = 3072.182
$ws.Range("K136").Value = 5314.531800000001
$ws.Range("L136").Value = 9216.545999999998
$ws.Range("M136").Value = -2764.531800000001
$ws.Range("N136").Value = -14316.546

$ws = $wb.Worksheets.Item("CUL")
# row 17
$ws.Range("H17").Value = 1242.7142
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 1366.5
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 4099.5
$ws.Range("M17").Value = -1331
$ws.Range("N17").Value = -4437.5
# row 137
$ws.Range("H137").Value = 40008348
$ws.Range("I137").Value = 2835.6428
$ws.Range("J137").Value = 90924450
$ws.Range("K137").Value = 8506.928400000001
$ws.Range("L137").Value = 272773350
$ws.Range("M137").Value = -3406.928400000001
$ws.Range("N137").Value = -272783550

$ws = $wb.Worksheets.Item("GSM")
# row 107
$ws.Range("H107").Value = 78414.766
$ws.Range("I107").Value = 80151.36
$ws.Range("K107").Value = 80151.36
$ws.Range("M107").Value = -78231.36
# row 122
$ws.Range("H122").Value = 1900.5333
$ws.Range("I122").Value = 1957.1428
$ws.Range("J122").Value = 1851
$ws.Range("K122").Value = 5871.428400000001
$ws.Range("L122").Value = 5553
$ws.Range("M122").Value = -3421.428400000001
$ws.Range("N122").Value = -10453
# row 126
$ws.Range("H126").Value = 5995.3335
$ws.Range("I126").Value = 11629.454
$ws.Range("J126").Value = 2121.875
$ws.Range("K126").Value = 34888.362
$ws.Range("L126").Value = 6365.625
$ws.Range("M126").Value = -32418.362
$ws.Range("N126").Value = -11305.625

$ws = $wb.Worksheets.Item("LTW")
# row 100
$ws.Range("H100").Value = 2936.75
$ws.Range("I100").Value = 2580
$ws.Range("J100").Value = 3531.3333
$ws.Range("K100").Value = 2580
$ws.Range("L100").Value = 3531.3333
$ws.Range("M100").Value = -2039
$ws.Range("N100").Value = -4613.3333
# row 132
$ws.Range("H132").Value = 3745.2703
$ws.Range("I132").Value = 3385.2144
$ws.Range("K132").Value = 10155.6432
$ws.Range("M132").Value = -7625.643199999999

$ws = $wb.Worksheets.Item("WVR")
# row 107
$ws.Range("H107").Value = 14286764
$ws.Range("I107").Value = 461.75
$ws.Range("J107").Value = 33335166
$ws.Range("K107").Value = 1385.25
$ws.Range("L107").Value = 100005498
$ws.Range("M107").Value = 534.75
$ws.Range("N107").Value = -100009338
# row 122
$ws.Range("H122").Value = 25170882
$ws.Range("I122").Value = 31093190
$ws.Range("J122").Value = 1075
$ws.Range("K122").Value = 93279570
$ws.Range("L122").Value = 3225
$ws.Range("M122").Value = -93277120
$ws.Range("N122").Value = -8125
# row 132
$ws.Range("H132").Value = 942.1404
$ws.Range("I132").Value = 781.5128
$ws.Range("J132").Value = 1290.1666
$ws.Range("K132").Value = 2344.5384
$ws.Range("L132").Value = 3870.4998
$ws.Range("M132").Value = 185.4616000000001
$ws.Range("N132").Value = -8930.4998
# row 141
$ws.Range("H141").Value = 44125
$ws.Range("J141").Value = 44125
$ws.Range("L141").Value = 44125
$ws.Range("N141").Value = -54485
